# Unaccounted-for staff window: duplicate the "Absent" status column (H)
# into two more trailing columns (I, J) and update the Absent-count
# column (F) to account for the 2 extra status columns now being summed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate column H into I and J (values + type, no new formatting) ---
$ws.Columns("H").Copy()
$ws.Columns("I").PasteSpecial(-4163)   # xlPasteValues
$ws.Columns("H").Copy()
$ws.Columns("J").PasteSpecial(-4163)   # xlPasteValues

# Column H is 11.65625 "width units" wide (best-fit); give I and J the same
# display width (closest representable ColumnWidth in this engine's model).
$ws.Columns("I").ColumnWidth = 10.83
$ws.Columns("J").ColumnWidth = 10.83

# --- Bump the Absent-count column (F) by 2 for every data row ---
$rows = @(4,5,6,7,9,10,11,12,14,15,16,17,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41)
foreach ($r in $rows) {
    $cell = $ws.Range("F" + $r)
    $cell.Value = $cell.Value() + 2
}

Write-Host "done"
